$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("ECs", "Hsp90aa1", "Fgfr3", "ECs", 3, 1, 72.31432066666666, 216.942962, 0.1293883843050027, 0.1293883843050027, 3, 1, 4.618552666666667, 13.855658, 0.78434648953826, 0.78434648953826, 333.9874985532218, 3005.887486978996, 0.1014853250166562, 0.1014853250166562),
    @("ECs", "Hsp90aa1", "Fgfr3", "FAPs", 3, 1, 72.31432066666666, 216.942962, 0.1293883843050027, 0.1293883843050027, 2, 0.6666666666666666, 0.6792986666666666, 2.037896, 0.1153620112191035, 0.1153620112191036, 49.12302160977244, 442.107194487952, 0.0149265042418154, 0.0149265042418154),
    @("ECs", "Hsp90aa1", "Fgfr3", "MuSCs", 3, 1, 72.31432066666666, 216.942962, 0.1293883843050027, 0.1293883843050027, 3, 1, 0.478937, 1.436811, 0.0813355572127976, 0.08133555721279762, 34.63400379713133, 311.706034174182, 0.01052387633431099, 0.01052387633431099),
    @("ECs", "Hsp90aa1", "Fgfr3", "Resolving-Mac", 3, 1, 72.31432066666666, 216.942962, 0.1293883843050027, 0.1293883843050027, 1, 0.3333333333333333, 0.1116203333333333, 0.334861, 0.01895594202983873, 0.01895594202983874, 8.071748577586888, 72.64573719828201, 0.002452678712220127, 0.002452678712220128),
    @("FAPs", "Hsp90aa1", "Fgfr3", "ECs", 3, 1, 154.0559436666667, 462.167831, 0.2756445674916043, 0.2756445674916043, 3, 1, 4.618552666666667, 13.855658, 0.78434648953826, 0.78434648953826, 711.5154894375332, 6403.639404937799, 0.2162008488723318, 0.2162008488723318),
    @("FAPs", "Hsp90aa1", "Fgfr3", "FAPs", 3, 1, 154.0559436666667, 462.167831, 0.2756445674916043, 0.2756445674916043, 2, 0.6666666666666666, 0.6792986666666666, 2.037896, 0.1153620112191035, 0.1153620112191036, 104.6499971248418, 941.8499741235761, 0.03179891168745139, 0.03179891168745141),
    @("FAPs", "Hsp90aa1", "Fgfr3", "MuSCs", 3, 1, 154.0559436666667, 462.167831, 0.2756445674916043, 0.2756445674916043, 3, 1, 0.478937, 1.436811, 0.0813355572127976, 0.08133555721279762, 73.78309149188235, 664.047823426941, 0.02241970448961023, 0.02241970448961024),
    @("FAPs", "Hsp90aa1", "Fgfr3", "Resolving-Mac", 3, 1, 154.0559436666667, 462.167831, 0.2756445674916043, 0.2756445674916043, 1, 0.3333333333333333, 0.1116203333333333, 0.334861, 0.01895594202983873, 0.01895594202983874, 17.19577578405456, 154.761982056491, 0.005225102442210821, 0.005225102442210822),
    @("Inflammatory-Mac", "Hsp90aa1", "Fgfr3", "ECs", 3, 1, 168.0546723333333, 504.1640170000001, 0.300691790058393, 0.300691790058393, 3, 1, 4.618552666666667, 13.855658, 0.78434648953826, 0.78434648953826, 776.1693550509096, 6985.524195458187, 0.235846549965276, 0.235846549965276),
    @("Inflammatory-Mac", "Hsp90aa1", "Fgfr3", "FAPs", 3, 1, 168.0546723333333, 504.1640170000001, 0.300691790058393, 0.300691790058393, 2, 0.6666666666666666, 0.6792986666666666, 2.037896, 0.1153620112191035, 0.1153620112191036, 114.1593148431369, 1027.433833588232, 0.03468840965820865, 0.03468840965820866),
    @("Inflammatory-Mac", "Hsp90aa1", "Fgfr3", "MuSCs", 3, 1, 168.0546723333333, 504.1640170000001, 0.300691790058393, 0.300691790058393, 3, 1, 0.478937, 1.436811, 0.0813355572127976, 0.08133555721279762, 80.48760060330967, 724.3884054297871, 0.02445693429371295, 0.02445693429371295),
    @("Inflammatory-Mac", "Hsp90aa1", "Fgfr3", "Resolving-Mac", 3, 1, 168.0546723333333, 504.1640170000001, 0.300691790058393, 0.300691790058393, 1, 0.3333333333333333, 0.1116203333333333, 0.334861, 0.01895594202983873, 0.01895594202983874, 18.75831854407078, 168.824866896637, 0.005699896141195335, 0.005699896141195337),
    @("MuSCs", "Hsp90aa1", "Fgfr3", "ECs", 3, 1, 65.818911, 197.456733, 0.11776647381174, 0.11776647381174, 3, 1, 4.618552666666667, 13.855658, 0.78434648953826, 0.78434648953826, 303.988106916146, 2735.892962245314, 0.09236972031953772, 0.09236972031953773),
    @("MuSCs", "Hsp90aa1", "Fgfr3", "FAPs", 3, 1, 65.818911, 197.456733, 0.11776647381174, 0.11776647381174, 2, 0.6666666666666666, 0.6792986666666666, 2.037896, 0.1153620112191035, 0.1153620112191036, 44.710698483752, 402.396286353768, 0.01358577727310421, 0.01358577727310422),
    @("MuSCs", "Hsp90aa1", "Fgfr3", "MuSCs", 3, 1, 65.818911, 197.456733, 0.11776647381174, 0.11776647381174, 3, 1, 0.478937, 1.436811, 0.0813355572127976, 0.08133555721279762, 31.523111777607, 283.708005998463, 0.009578601768464211, 0.009578601768464213),
    @("MuSCs", "Hsp90aa1", "Fgfr3", "Resolving-Mac", 3, 1, 65.818911, 197.456733, 0.11776647381174, 0.11776647381174, 1, 0.3333333333333333, 0.1116203333333333, 0.334861, 0.01895594202983873, 0.01895594202983874, 7.346728785457, 66.12055906911301, 0.002232374450633865, 0.002232374450633866),
    @("Resolving-Mac", "Hsp90aa1", "Fgfr3", "ECs", 3, 1, 98.64960366666666, 295.948811, 0.17650878433326, 0.1765087843332599, 3, 1, 4.618552666666667, 13.855658, 0.78434648953826, 0.78434648953826, 455.6183900802931, 4100.565510722638, 0.1384440453644583, 0.1384440453644583),
    @("Resolving-Mac", "Hsp90aa1", "Fgfr3", "FAPs", 3, 1, 98.64960366666666, 295.948811, 0.17650878433326, 0.1765087843332599, 2, 0.6666666666666666, 0.6792986666666666, 2.037896, 0.1153620112191035, 0.1153620112191036, 67.01254423796178, 603.1128981416559, 0.02036240835852386, 0.02036240835852387),
    @("Resolving-Mac", "Hsp90aa1", "Fgfr3", "MuSCs", 3, 1, 98.64960366666666, 295.948811, 0.17650878433326, 0.1765087843332599, 3, 1, 0.478937, 1.436811, 0.0813355572127976, 0.08133555721279762, 47.24694523130233, 425.222507081721, 0.01435644032669922, 0.01435644032669922),
    @("Resolving-Mac", "Hsp90aa1", "Fgfr3", "Resolving-Mac", 3, 1, 98.64960366666666, 295.948811, 0.17650878433326, 0.1765087843332599, 1, 0.3333333333333333, 0.1116203333333333, 0.334861, 0.01895594202983873, 0.01895594202983874, 11.01130164447456, 99.101714800271, 0.003345890283578583, 0.003345890283578583)
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $row = $data[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($r + 2, $c + 1).Value = $row[$c]
    }
}

